$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '321.91'
Set-TextValue $ws.Range('E2') '-3.17%'
Set-TextValue $ws.Range('D3') '42.90'
Set-TextValue $ws.Range('E3') '-5.79%'
Set-TextValue $ws.Range('D4') '5.206'
Set-TextValue $ws.Range('E4') '-6.06%'
Set-TextValue $ws.Range('D5') '0.08179'
Set-TextValue $ws.Range('E5') '-3.56%'
Set-TextValue $ws.Range('D6') '4.317'
Set-TextValue $ws.Range('E6') '-2.57%'
Set-TextValue $ws.Range('D7') '1.801'
Set-TextValue $ws.Range('E7') '-13.39%'
Set-TextValue $ws.Range('D8') '0.9493'
Set-TextValue $ws.Range('E8') '-4.02%'
Set-TextValue $ws.Range('D9') '0.1120'
Set-TextValue $ws.Range('E9') '-2.81%'
Set-TextValue $ws.Range('D10') '0.1873'
Set-TextValue $ws.Range('E10') '-3.20%'
Set-TextValue $ws.Range('D11') '0.09375'
Set-TextValue $ws.Range('E11') '-4.99%'
Set-TextValue $ws.Range('D12') '0.04622'
Set-TextValue $ws.Range('D13') '7.464'
Set-TextValue $ws.Range('E13') '-21.30%'
Set-TextValue $ws.Range('D14') '0.1057'
Set-TextValue $ws.Range('E14') '-0.27%'
Set-TextValue $ws.Range('D15') '0.001305'
Set-TextValue $ws.Range('E15') '2.02%'
Set-TextValue $ws.Range('D16') '0.005739'
Set-TextValue $ws.Range('E16') '-2.51%'
Set-TextValue $ws.Range('B17') 'LEO'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D17') '3.358'
Set-TextValue $ws.Range('E17') '-0.82%'
Set-TextValue $ws.Range('B18') 'BTSEToken'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D18') '2.529'
Set-TextValue $ws.Range('E18') '-0.47%'
Set-TextValue $ws.Range('B19') 'BitpandaEcosystemToken'
Set-TextValue $ws.Range('C19') 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue $ws.Range('D19') '0.3363'
Set-TextValue $ws.Range('E19') '0.27%'
Set-TextValue $ws.Range('B20') 'ProBitToken'
Set-TextValue $ws.Range('C20') 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue $ws.Range('D20') '0.1388'
Set-TextValue $ws.Range('E20') '0.39%'
Set-TextValue $ws.Range('B21') 'ZBToken'
Set-TextValue $ws.Range('C21') 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
Set-TextValue $ws.Range('D21') '0.2550'
Set-TextValue $ws.Range('E21') '-0.14%'
Set-TextValue $ws.Range('B22') 'CoinExToken'
Set-TextValue $ws.Range('C22') 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range('D22') '0.04148'
Set-TextValue $ws.Range('E22') '0.12%'
Set-TextValue $ws.Range('B23') 'BitKan'
Set-TextValue $ws.Range('C23') 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue $ws.Range('D23') '0.001250'
Set-TextValue $ws.Range('E23') '-3.97%'
Set-TextValue $ws.Range('B24') 'HotbitToken'
Set-TextValue $ws.Range('C24') 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue $ws.Range('D24') '0.004278'
Set-TextValue $ws.Range('E24') '-7.53%'
Set-TextValue $ws.Range('D25') '0.0001199'
Set-TextValue $ws.Range('E25') '-7.96%'
Set-TextValue $ws.Range('D26') '0.0002979'
Set-TextValue $ws.Range('E26') '-0.22%'
Set-TextValue $ws.Range('D38') '0.02669'
Set-TextValue $ws.Range('E38') '-2.14%'
Set-TextValue $ws.Range('D39') '0.05561'
Set-TextValue $ws.Range('E39') '-3.24%'
Set-TextValue $ws.Range('D40') '0.008140'
Set-TextValue $ws.Range('E40') '4.10%'
Set-TextValue $ws.Range('D41') '0.1400'
Set-TextValue $ws.Range('E41') '-2.52%'
Set-TextValue $ws.Range('D42') '0.006547'
Set-TextValue $ws.Range('E42') '-9.71%'
Set-TextValue $ws.Range('D43') '0.002109'
Set-TextValue $ws.Range('E43') '-0.84%'
Set-TextValue $ws.Range('D44') '0.007664'
Set-TextValue $ws.Range('E44') '-5.01%'
Set-TextValue $ws.Range('D45') '0.3204'
Set-TextValue $ws.Range('E45') '-9.92%'
Set-TextValue $ws.Range('D46') '0.00006740'
Set-TextValue $ws.Range('E47') '-0.22%'
Set-TextValue $ws.Range('D48') '0.003107'
Set-TextValue $ws.Range('E48') '-9.29%'
Set-TextValue $ws.Range('D49') '0.004100'
Set-TextValue $ws.Range('E49') '15.89%'
Set-TextValue $ws.Range('D50') '0.00002100'
Set-TextValue $ws.Range('E50') '-0.22%'
Set-TextValue $ws.Range('D51') '0.0002000'
Set-TextValue $ws.Range('E51') '-0.22%'
